$wb = $excel.ActiveWorkbook

# --- Metadata sheet updates ---
$meta = $wb.Worksheets.Item("Metadata")

# URL
$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/utilization-review"

# Version
$meta.Range("B3").Value = "8.0.0"

# Date
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"

# Publisher
$meta.Range("B9").Value = "LinuxForHealth Team"

# --- Elements sheet updates ---
$elem = $wb.Worksheets.Item("Elements")

# Row 2 is the base "Extension" element row; its Constraint(s) cell (AI2) used to
# duplicate the ele-1/ext-1 constraint text that belongs on "Extension.extension"
# (row 4, column AI). Clear the duplicate so it only appears once, on row 4.
$elem.Range("AI2").Value = ""

# Extension.url row's Fixed Value (Q5) mirrors the StructureDefinition URL shown
# on the Metadata sheet; keep it in sync with the new url.
$elem.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/utilization-review"
